# "input sheet burn with company name"
# This clears out the sample/demo data values that ship with the template
# so that the workbook becomes a blank input sheet ready to be filled in
# with a specific company's data. Clearing the cells that referenced the
# shared string "NA" also makes that shared-string entry unused, so it is
# dropped from the shared strings table on save and every subsequent
# shared-string index shifts down by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shareholding pattern table (rows 4-8, columns C:F) - drop sample years/values
$cellsToClear = @(
    "C4", "D4", "E4", "F4",
    "C5", "D5", "E5", "F5",
    "C6", "D6", "E6", "F6",
    "C7", "D7", "E7", "F7",
    "C8", "D8", "E8", "F8",

    # Board composition counts (rows 16-18, column C)
    "C16",
    "C17",
    "C18",

    # Indexed TSR / MD chart sample data (rows 38-42, columns B:D)
    "B38", "C38", "D38",
    "B39", "C39", "D39",
    "B40", "C40", "D40",
    "B41", "C41", "D41",
    "B42", "C42", "D42",

    # Executive compensation vs shareholder value table (rows 50-51, columns C:D)
    # C50/D50 held the now-removed "NA" shared string.
    "C50", "D50",
    "C51", "D51"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Dividend / EPS / Payout sample ratios (rows 28-29, columns C:D) reset to 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
